# "falta mergear filas en excel" - merge/reflow the response rows so the
# second block of answers (rows 5-7) continues the numbering of the first
# block (rows 2-4) instead of restarting, splitting the comma-joined
# "Det"/"Cont" values across the now length-6 sequence of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Lám 1, Rta 1)
$ws.Range("F2").Value = "C',FM"
$ws.Range("H2").Value = "?"

# Row 3 (Lám 1, Rta 2)
$ws.Range("F3").Value = "m"
$ws.Range("H3").Value = "?"
$ws.Range("I3").Value = "Fi"

# Row 4 (Lám 1, Rta 3)
$ws.Range("F4").Value = "M"
$ws.Range("H4").Value = "2"
$ws.Range("I4").Value = "H"
$ws.Range("J4").Value = "?"

# Row 5 (was Lám 2, Rta 1 -> renumbered to 4, continuing the single sequence)
$ws.Range("A5").Value = 4
$ws.Range("F5").Value = "?"
$ws.Range("H5").Value = "?"
$ws.Range("I5").Value = "Ad"
$ws.Range("J5").Value = "?"

# Row 6 (was Lám 2, Rta 2 -> renumbered to 5)
$ws.Range("A6").Value = 5
$ws.Range("F6").Value = "?"
$ws.Range("H6").Value = "?"
$ws.Range("I6").Value = "Hx"
$ws.Range("J6").Value = "?"

# Row 7 (was Lám 2, Rta 3 -> renumbered to 6)
$ws.Range("A7").Value = 6
$ws.Range("F7").Value = "?"
$ws.Range("H7").Value = "?"
$ws.Range("I7").Value = "Hd"
$ws.Range("J7").Value = "?"
